$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: wrap text (adds a new cellXf with alignment wrapText=1) ---
$ws.Range("A1:I1").WrapText = $true

# --- Data corrections (E5:E10, C11:C12, E11:E12, F12) ---
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 0

$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 0

$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1

# --- Formula changes in column I: simplified match check ---
$ws.Range("I2").Formula = "=IF(B2=G2, 1, 0)"
$ws.Range("I3").Formula = "=IF(B3=G3, 1, 0)"
$ws.Range("I4").Formula = "=IF(B4=G4, 1, 0)"
$ws.Range("I5").Formula = "=IF(B5=G5, 1, 0)"
$ws.Range("I6").Formula = "=IF(B6=G6, 1, 0)"
$ws.Range("I7").Formula = "=IF(B7=G7, 1, 0)"
$ws.Range("I8").Formula = "=IF(B8=G8, 1, 0)"
$ws.Range("I9").Formula = "=IF(B9=G9, 1, 0)"
$ws.Range("I10").Formula = "=IF(B10=G10, 1, 0)"
$ws.Range("I11").Formula = "=IF(B11=G11, 1, 0)"
$ws.Range("I12").Formula = "=IF(B12=G12, 1, 0)"

# --- Conditional formatting: 2-color scale on H2:H12 and I2:I12 ---
$csH = $ws.Range("H2:H12").FormatConditions.AddColorScale(2)
$csH.ColorScaleCriteria.Item(1).Type = 1
$csH.ColorScaleCriteria.Item(1).FormatColor.Color = 253
$csH.ColorScaleCriteria.Item(2).Type = 2
$csH.ColorScaleCriteria.Item(2).FormatColor.Color = 64768

$csI = $ws.Range("I2:I12").FormatConditions.AddColorScale(2)
$csI.ColorScaleCriteria.Item(1).Type = 1
$csI.ColorScaleCriteria.Item(1).FormatColor.Color = 253
$csI.ColorScaleCriteria.Item(2).Type = 2
$csI.ColorScaleCriteria.Item(2).FormatColor.Color = 64768
